# edit.ps1 - applies the FinalProjectPoster.pptx revision:
#  1. Abstract text box: drop the run-on "so" and expand the trailing
#     example list (construction/manmade sounds, singing and sounds of
#     nature); grow the text box to fit.
#  2. Conclusion text box: refocus the "given more time" paragraph on
#     features instead of decision-tree parameters, add "random forests"
#     to the list of models to test, and trim the closing sentence.
#
# (The master/layout footer "date" placeholders also show a refreshed
# 3/15/16 stamp in the canonical export, but that is PowerPoint's own
# whole-deck datetimeFigureOut re-cache on save, not a content edit --
# it is intentionally left alone here since the field can't be
# re-cached through the object model without destroying the <a:fld>.)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------
# 1. Abstract (shape "Subtitle 2", id 3)
# ---------------------------------------------------------------------
$abstract = $s.Shapes.Item(3)
$abstract.TextFrame.TextRange.Text = "    We built a machine learning system to identify whether a recording of a sound has a human scream in it. This task is important as it may potentially allow for computerized surveillance systems to replace human-centric ones. Since this in the context of a surveillance system, we trained and tested our system on sounds that would occur in a variety of surveillance situations, such as construction/manmade sounds, singing and sounds of nature (i.e. bird and other animal calls)."

# Grow the box's height to fit the longer text (width/position unchanged).
$abstract.Height = 617.5555

# ---------------------------------------------------------------------
# 2. Conclusion (shape "Subtitle 2", id 10)
# ---------------------------------------------------------------------
$conclusion = $s.Shapes.Item(10)
$tr = $conclusion.TextFrame.TextRange

$para1 = $tr.Paragraphs(1, 1)
$para1.Text = "    Given more time and resources, we would consider exploring more features such as relations between time segments (e.g. a scream should have no beat) and further cross validation."

$para2 = $tr.Paragraphs(2, 1)
$para2.Text = "    This also includes testing on different machine learning models, such as nearest neighbors, neural nets, random forests, etc. "

$para3 = $tr.Paragraphs(3, 1)
$para3.Text = "     Judging from the correct classification rate with 7-fold cross validation of 0.8234, the current model is successful."
